$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @("37-31=","51+41=","58+17=","14+9=","13+56=","76+14=","84-81=","26-23=","20+37=","70-47=","62+29=","40+23=","10+87=","97-90=","52+27=","66+32=","43+17=","73-47=","11+83=","48+14=","24+72=","89-87=","50+4=","43+7=","95-51=","69-5=","37+15=","93-21=","92-27=","85-67=","77-76=","13+33=","17+42=","4+52=","28-20=","86-72=","84-18=","29-18=","76-32=","78-10=","11+30=","69-3=","9+31=","57+4=","26+38=","16+54=","40+42=","43+0=","53+43=","70-22=","21+43=","76-48=","6+30=","50-38=","97+0=","15+23=","83-68=","46-38=","13+73=","77-32=","97-96=","68+20=","76+9=","5+82=","16+53=","99-59=","77-42=","17+70=","83-78=","87+6=","98-55=","25+45=","59-50=","90-10=","21+70=","91-50=","35+59=","74+13=","14-4=","33+57=","27+28=","91-85=","36+30=","3+86=","96-90=","79-0=","30+28=","61+4=","78+15=","88-37=","18-16=","59-2=","29+57=","32+19=","2+34=","30+3=","27-26=","17+57=","9+19=","98-77=")

$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
